# Auto-generated edit script applying updated market-price values
# to the FFXIV Leve profit tracker workbook (8 job sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 999326.4
$ws.Range("I9").Value = 1623459
$ws.Range("J9").Value = 714.2
$ws.Range("K9").Value = 1623459
$ws.Range("L9").Value = 714.2
$ws.Range("M9").Value = -1623290
$ws.Range("N9").Value = -1052.2

$ws.Range("H40").Value = 18218.924
$ws.Range("I40").Value = 5872
$ws.Range("J40").Value = 45999.5
$ws.Range("K40").Value = 5872
$ws.Range("L40").Value = 45999.5
$ws.Range("M40").Value = -5697
$ws.Range("N40").Value = -46349.5

$ws.Range("H51").Value = 40059.8
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H58").Value = 270.7143
$ws.Range("I58").Value = 270.7143
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 812.1428999999999
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = -662.1428999999999
$ws.Range("M58").ClearContents()

$ws.Range("H88").Value = 7124.75
$ws.Range("J88").Value = 7499.6665
$ws.Range("L88").Value = 7499.6665
$ws.Range("N88").Value = -8311.666499999999

$ws.Range("H91").Value = 7124.75
$ws.Range("J91").Value = 7499.6665
$ws.Range("L91").Value = 7499.6665
$ws.Range("N91").Value = -10307.6665

$ws.Range("H100").Value = 4017.3157
$ws.Range("I100").Value = 3097.4443
$ws.Range("J100").Value = 4845.2
$ws.Range("K100").Value = 3097.4443
$ws.Range("L100").Value = 4845.2
$ws.Range("M100").Value = -2556.4443
$ws.Range("N100").Value = -5927.2

$ws.Range("H112").Value = 4912.727
$ws.Range("I112").Value = 1433
$ws.Range("J112").Value = 6217.625
$ws.Range("K112").Value = 4299
$ws.Range("L112").Value = 18652.875
$ws.Range("M112").Value = -3191
$ws.Range("N112").Value = -20868.875

$ws.Range("H133").Value = 100290
$ws.Range("J133").Value = 100290
$ws.Range("L133").Value = 100290
$ws.Range("N133").Value = -110410

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6300.3125
$ws.Range("I45").Value = 8265.387000000001
$ws.Range("K45").Value = 8265.387000000001
$ws.Range("M45").Value = -7888.387000000001

$ws.Range("H63").Value = 1484.75
$ws.Range("I63").Value = 1484.75
$ws.Range("K63").Value = 1484.75
$ws.Range("M63").Value = -798.75

$ws.Range("H66").Value = 1484.75
$ws.Range("I66").Value = 1484.75
$ws.Range("K66").Value = 7423.75
$ws.Range("M66").Value = -3991.75

$ws.Range("H74").Value = 2825.4285
$ws.Range("I74").Value = 2754.1667
$ws.Range("J74").Value = 3253
$ws.Range("K74").Value = 2754.1667
$ws.Range("L74").Value = 3253
$ws.Range("M74").Value = -1880.1667
$ws.Range("N74").Value = -5001

$ws.Range("H76").Value = 70000
$ws.Range("J76").Value = 70000
$ws.Range("L76").Value = 70000
$ws.Range("N76").Value = -70676

$ws.Range("H77").Value = 2825.4285
$ws.Range("I77").Value = 2754.1667
$ws.Range("J77").Value = 3253
$ws.Range("K77").Value = 13770.8335
$ws.Range("L77").Value = 16265
$ws.Range("M77").Value = -9402.833500000001
$ws.Range("N77").Value = -25001

$ws.Range("H79").Value = 70000
$ws.Range("J79").Value = 70000
$ws.Range("L79").Value = 70000
$ws.Range("N79").Value = -72340

$ws.Range("H97").Value = 1900.6154
$ws.Range("I97").Value = 1840.8
$ws.Range("K97").Value = 1840.8
$ws.Range("M97").Value = -1344.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3741.3635
$ws.Range("I86").Value = 1609.6
$ws.Range("J86").Value = 5517.8335
$ws.Range("K86").Value = 1609.6
$ws.Range("L86").Value = 5517.8335
$ws.Range("M86").Value = -486.5999999999999
$ws.Range("N86").Value = -7763.8335

$ws.Range("H89").Value = 3741.3635
$ws.Range("I89").Value = 1609.6
$ws.Range("J89").Value = 5517.8335
$ws.Range("K89").Value = 8048
$ws.Range("L89").Value = 27589.1675
$ws.Range("M89").Value = -2432
$ws.Range("N89").Value = -38821.1675

$ws.Range("H94").Value = 1203.875
$ws.Range("I94").Value = 1165.1428
$ws.Range("K94").Value = 1165.1428
$ws.Range("M94").Value = -714.1428000000001

$ws.Range("H129").Value = 90000
$ws.Range("J129").Value = 90000
$ws.Range("L129").Value = 90000
$ws.Range("N129").Value = -100000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3038.7273
$ws.Range("J75").Value = 3182.5715
$ws.Range("L75").Value = 9547.7145
$ws.Range("N75").Value = -11543.7145

$ws.Range("H78").Value = 3038.7273
$ws.Range("J78").Value = 3182.5715
$ws.Range("L78").Value = 28643.1435
$ws.Range("N78").Value = -38627.1435

$ws.Range("H94").Value = 14692.286
$ws.Range("J94").Value = 14692.286
$ws.Range("L94").Value = 44076.858
$ws.Range("N94").Value = -45428.858

$ws.Range("H98").Value = 1791.3334
$ws.Range("J98").Value = 2089.5
$ws.Range("L98").Value = 6268.5
$ws.Range("N98").Value = -9264.5

$ws.Range("H131").Value = 2196.2856
$ws.Range("I131").Value = 1558.4445
$ws.Range("K131").Value = 4675.333500000001
$ws.Range("M131").Value = 364.6664999999994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("L26").ClearContents()

$ws.Range("H38").Value = 29998.5
$ws.Range("J38").Value = 29998.5
$ws.Range("L38").Value = 29998.5
$ws.Range("N38").Value = -30924.5

$ws.Range("H46").Value = 6000
$ws.Range("I46").Value = 6000
$ws.Range("K46").Value = 6000
$ws.Range("M46").Value = -5844

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("N50").Value = 0
$ws.Range("L50").ClearContents()

$ws.Range("H58").Value = 29511.75
$ws.Range("I58").Value = 29511.75
$ws.Range("K58").Value = 29511.75
$ws.Range("M58").Value = -29234.75

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("N69").Value = 0
$ws.Range("L69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("N72").Value = 0
$ws.Range("L72").ClearContents()

$ws.Range("H80").Value = 3952
$ws.Range("I80").Value = 3005
$ws.Range("J80").Value = 4899
$ws.Range("K80").Value = 3005
$ws.Range("L80").Value = 4899
$ws.Range("M80").Value = -2007
$ws.Range("N80").Value = -6895

$ws.Range("H83").Value = 3952
$ws.Range("I83").Value = 3005
$ws.Range("J83").Value = 4899
$ws.Range("K83").Value = 15025
$ws.Range("L83").Value = 24495
$ws.Range("M83").Value = -10033
$ws.Range("N83").Value = -34479

$ws.Range("H132").Value = 6769.4346
$ws.Range("I132").Value = 6435.294
$ws.Range("J132").Value = 7716.1665
$ws.Range("K132").Value = 19305.882
$ws.Range("L132").Value = 23148.4995
$ws.Range("M132").Value = -16775.882
$ws.Range("N132").Value = -28208.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3679.75
$ws.Range("J22").Value = 4165.3335
$ws.Range("L22").Value = 4165.3335
$ws.Range("N22").Value = -4755.3335

$ws.Range("H27").Value = 3679.75
$ws.Range("J27").Value = 4165.3335
$ws.Range("L27").Value = 4165.3335
$ws.Range("N27").Value = -4379.3335

$ws.Range("H55").Value = 462.8125
$ws.Range("I55").Value = 493.1
$ws.Range("J55").Value = 412.33334
$ws.Range("K55").Value = 493.1
$ws.Range("L55").Value = 412.33334
$ws.Range("M55").Value = -320.1
$ws.Range("N55").Value = -758.33334

$ws.Range("H100").Value = 4170.7896
$ws.Range("I100").Value = 3080.625
$ws.Range("J100").Value = 4963.636
$ws.Range("K100").Value = 3080.625
$ws.Range("L100").Value = 4963.636
$ws.Range("M100").Value = -2539.625
$ws.Range("N100").Value = -6045.636

$ws.Range("H132").Value = 5841.4546
$ws.Range("I132").Value = 4709.6665
$ws.Range("K132").Value = 14128.9995
$ws.Range("M132").Value = -11598.9995

$ws.Range("H136").Value = 1797.08
$ws.Range("I136").Value = 1704.2667
$ws.Range("K136").Value = 5112.800099999999
$ws.Range("M136").Value = -2562.800099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2200.7
$ws.Range("I100").Value = 2082.4
$ws.Range("K100").Value = 4164.8
$ws.Range("M100").Value = -3623.8

$ws.Range("H111").Value = 24000
$ws.Range("J111").Value = 24000
$ws.Range("L111").Value = 24000
$ws.Range("N111").Value = -32180

$ws.Range("H126").Value = 1593.6923
$ws.Range("I126").Value = 1593.6923
$ws.Range("K126").Value = 4781.0769
$ws.Range("M126").Value = -2311.0769

$ws.Range("H132").Value = 2504.15
$ws.Range("I132").Value = 2688.2727
$ws.Range("J132").Value = 1636.1428
$ws.Range("K132").Value = 8064.8181
$ws.Range("L132").Value = 4908.428400000001
$ws.Range("M132").Value = -5534.8181
$ws.Range("N132").Value = -9968.428400000001

$ws.Range("H139").Value = 121571.336
$ws.Range("J139").Value = 121571.336
$ws.Range("L139").Value = 121571.336
$ws.Range("N139").Value = -131851.336
